$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.836.37'
$ws.Range("E2").Value = '  +3.04%  '

# Row 3
$ws.Range("D3").Value = '1.869.02'
$ws.Range("E3").Value = '  +2.70%  '

# Row 4
$ws.Range("D4").Value = '''1.041'
$ws.Range("E4").Value = '  +3.25%  '

# Row 5
$ws.Range("D5").Value = '''324.96'
$ws.Range("E5").Value = '  +3.54%  '

# Row 6
$ws.Range("D6").Value = '''1.037'
$ws.Range("E6").Value = '  +2.92%  '

# Row 7
$ws.Range("D7").Value = '''0.4432'
$ws.Range("E7").Value = '  +2.84%  '

# Row 8
$ws.Range("E8").Value = '  +2.93%  '

# Row 9
$ws.Range("D9").Value = '''0.07482'
$ws.Range("E9").Value = '  +2.80%  '

# Row 10
$ws.Range("D10").Value = '''0.8874'
$ws.Range("E10").Value = '  +1.93%  '

# Row 11
$ws.Range("D11").Value = '''21.82'
$ws.Range("E11").Value = '  +2.35%  '

# Row 12
$ws.Range("D12").Value = '1.917.29'
$ws.Range("E12").Value = '  -9.13%  '

# Row 13
$ws.Range("D13").Value = '''5.563'
$ws.Range("E13").Value = '  +2.59%  '

# Row 14
$ws.Range("D14").Value = '''6.771'
$ws.Range("E14").Value = '  +1.98%  '

# Row 15
$ws.Range("D15").Value = '''0.07247'
$ws.Range("E15").Value = '  +4.11%  '

# Row 16
$ws.Range("E16").Value = '  +3.37%  '

# Row 17
$ws.Range("E17").Value = '  +2.81%  '

# Row 18
$ws.Range("D18").Value = '''0.000009175'
$ws.Range("E18").Value = '  +3.64%  '

# Row 19
$ws.Range("D19").Value = '''1.037'
$ws.Range("E19").Value = '  +2.91%  '

# Row 20
$ws.Range("D20").Value = '''15.59'
$ws.Range("E20").Value = '  +2.03%  '

# Row 21
$ws.Range("D21").Value = '27.862.31'
$ws.Range("E21").Value = '  +2.94%  '

# Row 22
$ws.Range("D22").Value = '''5.327'
$ws.Range("E22").Value = '  +2.18%  '

# Row 23
$ws.Range("E23").Value = '  +3.13%  '

# Row 24
$ws.Range("D24").Value = '''1.974'
$ws.Range("E24").Value = '  +4.56%  '

# Row 25
$ws.Range("D25").Value = '''158.95'
$ws.Range("E25").Value = '  +2.88%  '

# Row 26
$ws.Range("E26").Value = '  +2.62%  '

# Row 27
$ws.Range("D27").Value = '''1.992'
$ws.Range("E27").Value = '  +3.54%  '

# Row 28
$ws.Range("D28").Value = '''5.335'
$ws.Range("E28").Value = '  +1.71%  '

# Row 29
$ws.Range("D29").Value = '''117.87'
$ws.Range("E29").Value = '  +2.43%  '

# Row 30
$ws.Range("D30").Value = '''0.09116'
$ws.Range("E30").Value = '  +1.46%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''0.7802'
$ws.Range("E31").Value = '  +3.78%  '

# Row 32
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '''1.224'
$ws.Range("E32").Value = '  +3.25%  '

# Row 33
$ws.Range("D33").Value = '''3.096'
$ws.Range("E33").Value = '  +10.00%  '

# Row 34
$ws.Range("E34").Value = '  +3.61%  '

# Row 35
$ws.Range("D35").Value = '''1.039'
$ws.Range("E35").Value = '  +3.13%  '

# Row 36
$ws.Range("D36").Value = '''1.165'
$ws.Range("E36").Value = '  +3.40%  '

# Row 37
$ws.Range("D37").Value = '''0.02004'
$ws.Range("E37").Value = '  +3.84%  '

# Row 38
$ws.Range("E38").Value = '  +2.26%  '

# Row 39
$ws.Range("D39").Value = '''2.860'
$ws.Range("E39").Value = '  +4.04%  '

# Row 40
$ws.Range("D40").Value = '''0.5217'
$ws.Range("E40").Value = '  +1.67%  '

# Row 41
$ws.Range("D41").Value = '''0.1700'
$ws.Range("E41").Value = '  +2.73%  '

# Row 42
$ws.Range("D42").Value = '''6.912'
$ws.Range("E42").Value = '  +6.32%  '

# Row 43
$ws.Range("D43").Value = '''8.734'
$ws.Range("E43").Value = '  +4.69%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''109.99'
$ws.Range("E44").Value = '  +2.49%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''10.66'
$ws.Range("E45").Value = '  +2.50%  '

# Row 46
$ws.Range("E46").Value = '  +4.35%  '

# Row 47
$ws.Range("D47").Value = '''0.4725'
$ws.Range("E47").Value = '  +2.64%  '

# Row 48
$ws.Range("D48").Value = '''0.06458'
$ws.Range("E48").Value = '  +3.70%  '

# Row 49
$ws.Range("D49").Value = '''1.904'
$ws.Range("E49").Value = '  +2.76%  '

# Row 50
$ws.Range("D50").Value = '''40.07'
$ws.Range("E50").Value = '  +4.08%  '

# Row 51
$ws.Range("D51").Value = '''64.74'
$ws.Range("E51").Value = '  +1.48%  '
